$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.783482801709658
$ws.Range("D2").Value = 3.927836968117492
$ws.Range("E2").Value = 16.57547603266558
$ws.Range("F2").Value = 18.86300847045873
$ws.Range("G2").Value = 19.49503046825303
$ws.Range("H2").Value = 11.90160245558593
$ws.Range("K2").Value = 12.99559702752478
$ws.Range("O2").Value = 16.87606019416074
$ws.Range("B3").Value = 7.710554446013313
$ws.Range("D3").Value = 3.842835098026128
$ws.Range("E3").Value = 15.63098910006513
$ws.Range("F3").Value = 18.89903405192173
$ws.Range("G3").Value = 19.53550223441135
$ws.Range("H3").Value = 11.95889704294781
$ws.Range("K3").Value = 12.46118444306538
$ws.Range("O3").Value = 16.96421746401018
$ws.Range("B4").Value = 7.66726960309127
$ws.Range("D4").Value = 3.78918522582065
$ws.Range("E4").Value = 15.02568019808862
$ws.Range("F4").Value = 18.9286428217236
$ws.Range("G4").Value = 19.5717630050552
$ws.Range("H4").Value = 11.99667619489057
$ws.Range("K4").Value = 12.11948367569278
$ws.Range("O4").Value = 17.0238138743887
$ws.Range("B5").Value = 7.650025995242809
$ws.Range("D5").Value = 3.766974796303285
$ws.Range("E5").Value = 14.77289448169157
$ws.Range("F5").Value = 18.94258193792209
$ws.Range("G5").Value = 19.58938373818484
$ws.Range("H5").Value = 12.01272431768255
$ws.Range("K5").Value = 11.97696763192751
$ws.Range("O5").Value = 17.04946743714767
$ws.Range("B6").Value = 7.647187146867876
$ws.Range("D6").Value = 3.763266350665093
$ws.Range("E6").Value = 14.73055898290164
$ws.Range("F6").Value = 18.94500934839397
$ws.Range("G6").Value = 19.59248063239569
$ws.Range("H6").Value = 12.01542849598073
$ws.Range("K6").Value = 11.95310996963455
$ws.Range("O6").Value = 17.05380957315531
$ws.Range("B7").Value = 7.667035423556799
$ws.Range("D7").Value = 3.788887069965469
$ws.Range("E7").Value = 15.02229541881497
$ws.Range("F7").Value = 18.9288232385675
$ws.Range("G7").Value = 19.57198916411369
$ws.Range("H7").Value = 11.99688998359308
$ws.Range("K7").Value = 12.11757469362879
$ws.Range("O7").Value = 17.02415431992766
$ws.Range("B8").Value = 7.758039216721575
$ws.Range("D8").Value = 3.898842101859296
$ws.Range("E8").Value = 16.2552312509162
$ws.Range("F8").Value = 18.87387022229164
$ws.Range("G8").Value = 19.50660259741881
$ws.Range("H8").Value = 11.92081739058325
$ws.Range("K8").Value = 12.81422711096848
$ws.Range("O8").Value = 16.90531767221484
$ws.Range("B9").Value = 7.947446087135011
$ws.Range("D9").Value = 4.102111626430481
$ws.Range("E9").Value = 18.57332170634794
$ws.Range("F9").Value = 18.82590900322502
$ws.Range("G9").Value = 19.46988510510188
$ws.Range("H9").Value = 11.79231418829604
$ws.Range("K9").Value = 14.067396117934
$ws.Range("O9").Value = 16.71598489612077
$ws.Range("B10").Value = 8.092034607590193
$ws.Range("D10").Value = 4.242959468758587
$ws.Range("E10").Value = 20.22245640928284
$ws.Range("F10").Value = 18.82755751252365
$ws.Range("G10").Value = 19.49976252700736
$ws.Range("H10").Value = 11.71056888348474
$ws.Range("K10").Value = 14.913304554876
$ws.Range("O10").Value = 16.60396359196894
$ws.Range("B11").Value = 8.158720072608922
$ws.Range("D11").Value = 4.305009166177397
$ws.Range("E11").Value = 20.93019459660107
$ws.Range("F11").Value = 18.83637817569645
$ws.Range("G11").Value = 19.52584918046143
$ws.Range("H11").Value = 11.67614505505422
$ws.Range("K11").Value = 15.28092631359639
$ws.Range("O11").Value = 16.5589745391402
$ws.Range("B12").Value = 8.18408033117853
$ws.Range("D12").Value = 4.328201259223694
$ws.Range("E12").Value = 21.19212706694419
$ws.Range("F12").Value = 18.8408814567659
$ws.Range("G12").Value = 19.53753050473295
$ws.Range("H12").Value = 11.66350815499261
$ws.Range("K12").Value = 15.41759938368531
$ws.Range("O12").Value = 16.54280425090584
$ws.Range("B13").Value = 8.178614113369905
$ws.Range("D13").Value = 4.323220192076574
$ws.Range("E13").Value = 21.13598477678402
$ws.Range("F13").Value = 18.83985984284857
$ws.Range("G13").Value = 19.53493448646974
$ws.Range("H13").Value = 11.66621198439661
$ws.Range("K13").Value = 15.3882781796932
$ws.Range("O13").Value = 16.54624817257637
$ws.Range("B14").Value = 8.160804445673659
$ws.Range("D14").Value = 4.306923377677229
$ws.Range("E14").Value = 20.95186545313132
$ws.Range("F14").Value = 18.83672534457738
$ws.Range("G14").Value = 19.52677405764848
$ws.Range("H14").Value = 11.67509741182964
$ws.Range("K14").Value = 15.29222168143801
$ws.Range("O14").Value = 16.55762679898997
$ws.Range("B15").Value = 8.149908879094562
$ws.Range("D15").Value = 4.296901022688565
$ws.Range("E15").Value = 20.83829720037266
$ws.Range("F15").Value = 18.83495688976038
$ws.Range("G15").Value = 19.52201045662056
$ws.Range("H15").Value = 11.68059195189501
$ws.Range("K15").Value = 15.23305209177127
$ws.Range("O15").Value = 16.56470955008205
$ws.Range("B16").Value = 8.087692915978643
$ws.Range("D16").Value = 4.238862457988723
$ws.Range("E16").Value = 20.17535247181081
$ws.Range("F16").Value = 18.82714371771516
$ws.Range("G16").Value = 19.49830975252567
$ws.Range("H16").Value = 11.71287427042452
$ws.Range("K16").Value = 14.88892789603399
$ws.Range("O16").Value = 16.60702451786848
$ws.Range("B17").Value = 8.049742392347897
$ws.Range("D17").Value = 4.202729359443217
$ws.Range("E17").Value = 19.75780571246779
$ws.Range("F17").Value = 18.8244200035238
$ws.Range("G17").Value = 19.48697610001364
$ws.Range("H17").Value = 11.73338696860096
$ws.Range("K17").Value = 14.67336668397417
$ws.Range("O17").Value = 16.63451789977459
$ws.Range("B18").Value = 8.028001427017495
$ws.Range("D18").Value = 4.181756931850781
$ws.Range("E18").Value = 19.51364521165783
$ws.Range("F18").Value = 18.82361300154592
$ws.Range("G18").Value = 19.48163329230585
$ws.Range("H18").Value = 11.74544526532587
$ws.Range("K18").Value = 14.54776714125132
$ws.Range("O18").Value = 16.65089284357887
$ws.Range("B19").Value = 8.020655987736426
$ws.Range("D19").Value = 4.174623902391036
$ws.Range("E19").Value = 19.43028887148446
$ws.Range("F19").Value = 18.8234701247004
$ws.Range("G19").Value = 19.48002603786005
$ws.Range("H19").Value = 11.74957259910874
$ws.Range("K19").Value = 14.50496619892317
$ws.Range("O19").Value = 16.65653331838995
$ws.Range("B20").Value = 8.053773436262764
$ws.Range("D20").Value = 4.20659551929297
$ws.Range("E20").Value = 19.80266788109809
$ws.Range("F20").Value = 18.82463130774111
$ws.Range("G20").Value = 19.48806081697508
$ws.Range("H20").Value = 11.73117644439457
$ws.Range("K20").Value = 14.69648114228991
$ws.Range("O20").Value = 16.63153301587342
$ws.Range("B21").Value = 8.166032831816938
$ws.Range("D21").Value = 4.311718524199269
$ws.Range("E21").Value = 21.00611030668861
$ws.Range("F21").Value = 18.83761444467027
$ws.Range("G21").Value = 19.52912201598391
$ws.Range("H21").Value = 11.6724767155978
$ws.Range("K21").Value = 15.32050513909528
$ws.Range("O21").Value = 16.55426105522469
$ws.Range("B22").Value = 8.240019237801564
$ws.Range("D22").Value = 4.378640308679356
$ws.Range("E22").Value = 21.75724889708506
$ws.Range("F22").Value = 18.85287907342289
$ws.Range("G22").Value = 19.56646630885645
$ws.Range("H22").Value = 11.63643740056553
$ws.Range("K22").Value = 15.71352802281823
$ws.Range("O22").Value = 16.50881102057602
$ws.Range("B23").Value = 8.200482331225208
$ws.Range("D23").Value = 4.343090275585697
$ws.Range("E23").Value = 21.35957848548738
$ws.Range("F23").Value = 18.84411131966831
$ws.Range("G23").Value = 19.54557249857464
$ws.Range("H23").Value = 11.65545910968157
$ws.Range("K23").Value = 15.5051388357382
$ws.Range("O23").Value = 16.53260385615031
$ws.Range("B24").Value = 8.051950758783608
$ws.Range("D24").Value = 4.204848247272122
$ws.Range("E24").Value = 19.78239849112612
$ws.Range("F24").Value = 18.82453341334368
$ws.Range("G24").Value = 19.48756676286772
$ws.Range("H24").Value = 11.73217499614131
$ws.Range("K24").Value = 14.68603629398421
$ws.Range("O24").Value = 16.63288071113758
$ws.Range("B25").Value = 7.895161252534678
$ws.Range("D25").Value = 4.048548186874972
$ws.Range("E25").Value = 17.92813771648611
$ws.Range("F25").Value = 18.83243168341431
$ws.Range("G25").Value = 19.46989703513598
$ws.Range("H25").Value = 11.82485802609753
$ws.Range("K25").Value = 13.74114521941061
$ws.Range("O25").Value = 16.76247993368605
